$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "SIS-ESTIVAL"
$ws.Range("C1").Value = "SIS-INVERNAL"
$ws.Range("A2").Value = "Mean"
$ws.Range("A8").Value = "Máx"
